# Adding test case: Check user can change Personal details

$wb = $excel.ActiveWorkbook

# Add the new "PersonalDetails" worksheet right after the existing "SignIn" sheet;
# Excel activates/selects a newly inserted sheet automatically.
$signInSheet = $wb.Worksheets.Item(1)
$newSheet = $wb.Worksheets.Add($null, $signInSheet)
$newSheet.Name = "PersonalDetails"

# Column headers (first/last name) entered first.
$newSheet.Range("A1").Value = "firstName"
$newSheet.Range("B1").Value = "lastName"

# Name data, row by row.
$newSheet.Range("A2").Value = "Anthony"
$newSheet.Range("B2").Value = "Regis"

$newSheet.Range("A3").Value = "Antonio"
$newSheet.Range("B3").Value = "Banderas"

$newSheet.Range("A4").Value = "Tony"
$newSheet.Range("B4").Value = "Stark"

$newSheet.Range("B5").Value = "Bonito"

$newSheet.Range("B6").Value = "King"
$newSheet.Range("A6").Value = "Antoine"

# Remaining headers.
$newSheet.Range("E1").Value = "faxNumber"
$newSheet.Range("D1").Value = "telephoneNumber"
$newSheet.Range("C1").Value = "email"

# Final fix-up of the first name in row 5.
$newSheet.Range("A5").Value = "Chono"

# Email formulas derived from the first name.
$newSheet.Range("C2").Formula = '=CONCATENATE(LOWER(A2), "@codifyme.co.nz")'
$newSheet.Range("C3").Formula = '=CONCATENATE(LOWER(A3), "@codifyme.co.nz")'
$newSheet.Range("C4").Formula = '=CONCATENATE(LOWER(A4), "@codifyme.co.nz")'
$newSheet.Range("C5").Formula = '=CONCATENATE(LOWER(A5), "@codifyme.co.nz")'
$newSheet.Range("C6").Formula = '=CONCATENATE(LOWER(A6), "@codifyme.co.nz")'

# Phone / fax numbers.
$newSheet.Range("D2").Value = 64221328444
$newSheet.Range("E2").Value = 6495551234

$newSheet.Range("D3").Value = 64221328445
$newSheet.Range("E3").Value = 6435555678

$newSheet.Range("D4").Value = 64221328446
$newSheet.Range("E4").Value = 6475559876

$newSheet.Range("D5").Value = 64221328447
$newSheet.Range("E5").Value = 6445554321

$newSheet.Range("D6").Value = 64221328448
$newSheet.Range("E6").Value = 6465558765

# Column widths to match bestFit sizing (values chosen so the engine's
# pixel-quantised stored width lands as close as possible to the target).
$newSheet.Columns.Item(1).ColumnWidth = 9.0
$newSheet.Columns.Item(2).ColumnWidth = 8.666666666666666
$newSheet.Columns.Item(3).ColumnWidth = 29.166666666666664
$newSheet.Columns.Item(4).ColumnWidth = 17.0
$newSheet.Columns.Item(5).ColumnWidth = 10.166666666666666

# Selection on the new, now-active sheet.
$newSheet.Range("C9").Select()
